# Updated capital structure database
# Refresh the Saudi Arabia Insurance (Life) dataset: revised growth/margin/
# valuation figures for the existing companies, a company-name correction for
# row 4 (now Saudi Enaya Cooperative Insurance Company) and a brand-new row 5
# entry for Al Sagr Cooperative Insurance Company. A couple of now-unused
# metric cells (F2/F3 expected EPS growth, T2/T3 buybacks-of-cash-returned)
# are cleared entirely rather than zeroed, matching the upstream schema change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 ----
$ws.Range("D2").Value = 0.113
$ws.Range("E2").Value = 0.0554
$ws.Range("F2").ClearContents()
$ws.Range("G2").Value = 0.05165445404659415
$ws.Range("H2").Value = 0.05165445404659415
$ws.Range("I2").Value = 0.06591528932408898
$ws.Range("J2").Value = 0.06222837507171944
$ws.Range("K2").Value = 158.1
$ws.Range("L2").Value = 0.05195018565373114
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("T2").ClearContents()
$ws.Range("U2").Value = 193.7
$ws.Range("V2").Value = 0.04614211867838681
$ws.Range("W2").Value = -0.223155929038282
$ws.Range("X2").Value = 0.0614077660836822
$ws.Range("Y2").Value = -0.2845636951219642
$ws.Range("Z2").Value = 6.064766839378239
$ws.Range("AA2").Value = -0.4019933554817276
$ws.Range("AB2").Value = 0.0614077660836822
$ws.Range("AC2").Value = -0.4634260165547759
$ws.Range("AD2").Value = 0.585
$ws.Range("AF2").Value = 0.585
$ws.Range("AG2").Value = -193.115
$ws.Range("AH2").Value = 0.000139335974762325
$ws.Range("AI2").Value = 0.000516882623466471
$ws.Range("AJ2").Value = -0.04822106555033541
$ws.Range("AK2").Value = -0.2058608761466178
$ws.Range("AN2").Value = 0.002845330739299611
$ws.Range("AP2").Value = -0.9392752918287937

# ---- Row 3 ----
$ws.Range("D3").Value = 0.113
$ws.Range("E3").Value = 0.0554
$ws.Range("F3").ClearContents()
$ws.Range("G3").Value = 0.06709784525991767
$ws.Range("H3").Value = 0.06709784525991767
$ws.Range("I3").Value = 0.08120914467540552
$ws.Range("J3").Value = 0.0675820565143707
$ws.Range("K3").Value = 195.4
$ws.Range("L3").Value = 0.0675820565143707
$ws.Range("M3").Value = -0
$ws.Range("N3").Value = -0
$ws.Range("O3").Value = -0
$ws.Range("P3").Value = -0
$ws.Range("Q3").Value = -0
$ws.Range("R3").Value = -0
$ws.Range("T3").ClearContents()
$ws.Range("U3").Value = 95.09999999999999
$ws.Range("V3").Value = 0.02442281517244922
$ws.Range("W3").Value = 0.2494255808016339
$ws.Range("X3").Value = 0.0614077660836822
$ws.Range("Y3").Value = 0.1880178147179517
$ws.Range("Z3").Value = 6.578612059158136
$ws.Range("AA3").Value = 0.4445961319681457
$ws.Range("AB3").Value = 0.0614077660836822
$ws.Range("AC3").Value = 0.3831883658844635
$ws.Range("AD3").Value = 0
$ws.Range("AF3").Value = 0
$ws.Range("AG3").Value = -95.09999999999999
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = -0.02503422133305254
$ws.Range("AK3").Value = -0.1035496515679443
$ws.Range("AN3").Value = 0
$ws.Range("AP3").Value = -0.398074508162411

# ---- Row 4 ----
$ws.Range("B4").Value = "Saudi Enaya Cooperative Insurance Company (SASE:8311)"
$ws.Range("D4").Value = 0.76
$ws.Range("G4").Value = -0.581140350877193
$ws.Range("H4").Value = -0.581140350877193
$ws.Range("I4").Value = -0.2653508771929824
$ws.Range("J4").Value = -0.2653508771929824
$ws.Range("K4").Value = -13.4
$ws.Range("L4").Value = -0.293859649122807
$ws.Range("U4").Value = 14.4
$ws.Range("V4").Value = 0.1044234952864394
$ws.Range("W4").Value = -0.2857142857142858
$ws.Range("X4").Value = 0.06158460752332806
$ws.Range("Y4").Value = -0.3472988932376138
$ws.Range("Z4").Value = 1.514950166112957
$ws.Range("AA4").Value = -0.4019933554817276
$ws.Range("AB4").Value = 0.06143266107304832
$ws.Range("AC4").Value = -0.4634260165547759
$ws.Range("AD4").Value = 0.585
$ws.Range("AF4").Value = 0.585
$ws.Range("AG4").Value = -13.815
$ws.Range("AH4").Value = 0.004224284218507419
$ws.Range("AI4").Value = 0.0172134765337649
$ws.Range("AJ4").Value = -0.1113349719950034
$ws.Range("AK4").Value = -0.7053867755935667
$ws.Range("AN4").Value = -0.05043103448275862
$ws.Range("AP4").Value = 1.190948275862069

# ---- Row 5 ----
$ws.Range("B5").Value = "Al Sagr Cooperative Insurance Company (SASE:8180)"
$ws.Range("D5").Value = 0.04219999999999999
$ws.Range("G5").Value = -0.09680451127819549
$ws.Range("H5").Value = -0.09680451127819549
$ws.Range("I5").Value = -0.2077067669172932
$ws.Range("J5").Value = -0.2077067669172932
$ws.Range("K5").Value = -23.9
$ws.Range("L5").Value = -0.2246240601503759
$ws.Range("U5").Value = 84.2
$ws.Range("V5").Value = 0.5069235400361228
$ws.Range("W5").Value = -0.223155929038282
$ws.Range("X5").Value = 0.0614077660836822
$ws.Range("Y5").Value = -0.2845636951219642
$ws.Range("Z5").Value = 3.304347826086958
$ws.Range("AA5").Value = -0.6863354037267083
$ws.Range("AB5").Value = 0.0614077660836822
$ws.Range("AC5").Value = -0.7477431698103905
$ws.Range("AG5").Value = -84.2
$ws.Range("AJ5").Value = -1.028083028083028
$ws.Range("AK5").Value = -842.0000000000479
$ws.Range("AP5").Value = 3.880184331797235

